# Update stock counters for several Lithuania 2€ commemorative coin rows
# (№10265 from 25.03.2024, https://2eurostore.ru/)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2€")

$ws.Range("H11").Value = 1
$ws.Range("H12").Value = 1
$ws.Range("H13").Value = 1
$ws.Range("H14").Value = 1
$ws.Range("H17").Value = 1
$ws.Range("H18").Value = 1
$ws.Range("H19").Value = 1

# Move the selection to reflect where the user left off editing
$ws.Activate()
$ws.Range("G23").Select()
